# Add a new admin user row (Rishan) to the admin_users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rishan"
$ws.Range("B2").Value = "rishangupta857@gmail.com"
$ws.Range("C2").Value = 'scrypt:32768:8:1$MDh6i4G7pSXpI7at$46bd5c5cac3202059b697604e75ff03acc9ed1db328f6e12ec42079df876c77fb1429d8c7999e80d4da398501f50dfdf90764d8934a129b84fc28a836bfad6b4'
